$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds date serial 45183 for rows 2-33; update to 45184.
$range = $ws.Range("C2:C33")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
